$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 92
$ws.Range("C3").Value = 7
$ws.Range("C4").Value = 64
$ws.Range("C5").Value = 6
$ws.Range("C6").Value = 28
$ws.Range("C7").Value = 9
